$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style from an existing header cell (AC1) onto the new header cells
# so "Wins"/"Losses"/"Ties" match the look of the rest of the header row.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Add new header cells for the team record columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team's win/loss/tie record for every player row (rows 2 through 56)
$lastRow = 56
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 86   # AD column -> Wins
    $ws.Cells.Item($r, 31).Value = 76   # AE column -> Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF column -> Ties
}
